# Add blank separator rows into the test data import sheet:
#   - one blank row after "Test_5" (before "Test_6")
#   - one blank row after "Test_8" (before "Test_9")
# This pushes the remaining rows down and grows the used range from
# A1:F11 to A1:F13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row before current row 7 (row that holds "Test_6").
# Existing rows 7..11 shift down to 8..12.
$ws.Rows.Item(7).Insert()

# Insert a blank row before current row 11 (row that now holds "Test_9",
# after the first insertion shifted it down from 10 to 11).
# Existing rows 11..12 shift down to 12..13.
$ws.Rows.Item(11).Insert()

# Mirror the recorded selection state in the target sheet.
$ws.Range("C20").Select()
